$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds numeric-looking text (e.g. "229.17", "41.679.32")
# that must stay stored as text, exactly like the original cells. Force a
# text number format while writing the values, then clear the formatting
# again so the cells end up with no explicit style (matching the source).
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '41.679.32'
$ws.Range("E2").Value = '  +5.51%  '
$ws.Range("D3").Value = '2.234.96'
$ws.Range("E3").Value = '  +3.42%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = '229.17'
$ws.Range("E5").Value = '  +0.50%  '
$ws.Range("D6").Value = '0.625'
$ws.Range("E6").Value = '  -1.39%  '
$ws.Range("D7").Value = '61.82'
$ws.Range("E7").Value = '  -2.52%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").Value = '0.403'
$ws.Range("E9").Value = '  +2.30%  '
$ws.Range("D10").Value = '59.01'
$ws.Range("E10").Value = '  +1.57%  '
$ws.Range("D11").Value = '0.0880'
$ws.Range("E11").Value = '  +3.35%  '
$ws.Range("E12").Value = '  +0.26%  '
$ws.Range("D13").Value = '2.566.73'
$ws.Range("E13").Value = '  +3.34%  '
$ws.Range("D14").Value = '15.71'
$ws.Range("E14").Value = '  -1.74%  '
$ws.Range("D15").Value = '22.13'
$ws.Range("E15").Value = '  +0.53%  '
$ws.Range("D16").Value = '0.801'
$ws.Range("E16").Value = '  -0.84%  '
$ws.Range("D17").Value = '5.59'
$ws.Range("E17").Value = '  +1.75%  '
$ws.Range("D18").Value = '2.222.61'
$ws.Range("E18").Value = '  +2.99%  '
$ws.Range("D19").Value = '41.567.35'
$ws.Range("E19").Value = '  +5.13%  '
$ws.Range("D20").Value = '73.40'
$ws.Range("E20").Value = '  +1.80%  '
$ws.Range("E21").Value = '  +6.77%  '
$ws.Range("D22").Value = '6.04'
$ws.Range("E22").Value = '  -2.23%  '
$ws.Range("D23").Value = '248.33'
$ws.Range("E23").Value = '  +8.10%  '
$ws.Range("D25").Value = '2.41'
$ws.Range("E25").Value = '  +3.57%  '
$ws.Range("D26").Value = '2.37'
$ws.Range("E26").Value = '  +0.99%  '
$ws.Range("D27").Value = '9.58'
$ws.Range("E27").Value = '  -0.85%  '
$ws.Range("B28").Value = 'Kaspa'
$ws.Range("C28").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D28").Value = '0.144'
$ws.Range("E28").Value = '  +4.03%  '
$ws.Range("B29").Value = 'Monero'
$ws.Range("C29").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D29").Value = '168.72'
$ws.Range("E29").Value = '  -1.90%  '
$ws.Range("D30").Value = '20.15'
$ws.Range("E30").Value = '  +1.50%  '
$ws.Range("E31").Value = '  +0.85%  '
$ws.Range("D32").Value = '2.80'
$ws.Range("E32").Value = '  +5.71%  '
$ws.Range("D33").Value = '0.123'
$ws.Range("E33").Value = '  -0.22%  '
$ws.Range("D34").Value = '5.00'
$ws.Range("E34").Value = '  +6.60%  '
$ws.Range("D35").Value = '4.64'
$ws.Range("E35").Value = '  +1.19%  '
$ws.Range("D36").Value = '0.0627'
$ws.Range("E36").Value = '  +1.60%  '
$ws.Range("D37").Value = '3.79'
$ws.Range("E37").Value = '  +4.15%  '
$ws.Range("D38").Value = '6.70'
$ws.Range("E38").Value = '  -4.55%  '
$ws.Range("D39").Value = '2.39'
$ws.Range("E39").Value = '  -1.66%  '
$ws.Range("D40").Value = '1.00'
$ws.Range("E40").Value = '  -0.09%  '
$ws.Range("D41").Value = '0.000237'
$ws.Range("E41").Value = '  +31.12%  '
$ws.Range("D42").Value = '4.91'
$ws.Range("E42").Value = '  +6.25%  '
$ws.Range("D43").Value = '0.0237'
$ws.Range("E43").Value = '  +4.27%  '
$ws.Range("D44").Value = '8.59'
$ws.Range("E44").Value = '  +11.10%  '
$ws.Range("D45").Value = '100.19'
$ws.Range("E45").Value = '  -2.10%  '
$ws.Range("D46").Value = '0.0964'
$ws.Range("E46").Value = '  +4.37%  '
$ws.Range("D47").Value = '1.487.65'
$ws.Range("E47").Value = '  -2.43%  '
$ws.Range("B48").Value = 'TrustWalletToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D48").Value = '1.19'
$ws.Range("E48").Value = '  -0.82%  '
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").Value = '1.19'
$ws.Range("E49").Value = '  -5.38%  '
$ws.Range("D50").Value = '2.77'
$ws.Range("E50").Value = '  -1.14%  '
$ws.Range("E51").Value = '  -1.70%  '

$priceRange.ClearFormats()
